# Got rid of the Environment Slide
#
# The "Environment" slide is removed entirely. The following
# "Reward System" slide (which shifts up to take its place) gets
# its first bullet point text updated to mention seeking higher rewards.

$p = $ppt.ActivePresentation

# Locate and delete the "Environment" slide (walk backwards so deleting
# doesn't disturb indices still to be examined).
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    $titleShape = $slide.Shapes.Item(1)
    if ($titleShape.HasTextFrame -and $titleShape.TextFrame.TextRange.Text -eq "Environment") {
        $slide.Delete()
    }
}

# Update the "Reward System" slide's first bullet.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $titleShape = $slide.Shapes.Item(1)
    if ($titleShape.HasTextFrame -and $titleShape.TextFrame.TextRange.Text -eq "Reward System") {
        $body = $slide.Shapes.Item(2).TextFrame.TextRange
        $firstRun = $body.Paragraphs(1).Runs(1)
        $firstRun.Text = "Used to let the bot know how it is doing and adjust to seek higher rewards"
    }
}
